# Actualización automática hashcode dom ene 13 01:39:43 CET 2019
# Update the hashcode values (column B) for a set of rows identified by
# their cell coordinates on the active ("hashcode.csv") worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "49fd61a0f117a3b9ae7042347894fb62"
$ws.Range("B17").Value = "6346ebf65402181973efc4282efd6afa"
$ws.Range("B44").Value = "801441973795329981791b8136476d51"
$ws.Range("B74").Value = "81b8198663d8342ceb3b8c0f92fab114"
$ws.Range("B89").Value = "1616b5e7f8bed5b4d7aed86321c8e87e"
$ws.Range("B99").Value = "934acdaaaa0b3be31f1a4c83585356c0"
$ws.Range("B110").Value = "391e31b1a8bd2400f63b4fbdf2ed30bd"
$ws.Range("B136").Value = "145f6cdd9e574970a49058607a4c57c6"
$ws.Range("B159").Value = "dbfc21f7e94c2499a7e91e097f364003"
$ws.Range("B161").Value = "43b27c02768b9c7c3fa9e56208ca190b"
$ws.Range("B168").Value = "a1b0e2550e24d1d6623b2a13cb8c46cb"
$ws.Range("B169").Value = "d8e2d3b430620fbcc36650018a5d213d"
$ws.Range("B227").Value = "79d7ac27c02b8ee4b146a8ebaf9cdac1"
$ws.Range("B232").Value = "ae22bcdb5a3d16e8e1bb7667b80435a8"
$ws.Range("B278").Value = "c471259a9ae3506bba77c0b291834b56"
$ws.Range("B281").Value = "7f6ab24a2600337270ff3e0396ae3efd"
$ws.Range("B302").Value = "0f1ef506e706195dbd93c49065f789b1"
$ws.Range("B339").Value = "4355b8ccd9f3d91560badc347230afcd"
$ws.Range("B345").Value = "d1f32890b74c9e8aba42588b693f86cc"
$ws.Range("B419").Value = "afba4ee92bb44bede48ddf483ac24705"
$ws.Range("B460").Value = "ef3bb11c9a11290215fab20c3653025e"
$ws.Range("B478").Value = "19b25a4ce25f6f97839a85d363ab88b0"
$ws.Range("B500").Value = "90638a5840cb2ea45547ac598d99705e"
$ws.Range("B501").Value = "10add39a694426657601535a2ecb2c04"
$ws.Range("B517").Value = "d58681c86cbed19c395aab18d70338ab"
$ws.Range("B550").Value = "8aab137630c87b0adee966d8555f7e13"
$ws.Range("B616").Value = "078638d89707ef761041c1aa1f6eb798"
$ws.Range("B627").Value = "0225aa8685f6b6a513936ce0d53587e9"
$ws.Range("B715").Value = "241d6411b78f5716839b34f023ec7a59"
$ws.Range("B768").Value = "b45c8bde2cac9396d620eb045d985164"
$ws.Range("B816").Value = "dc3ff660a48a009b2c263afaeeb131db"
$ws.Range("B819").Value = "ddcecae74f700d34aeb688e4eafe9966"
$ws.Range("B825").Value = "ee144aaf330dcd969107a5068c1f5d28"
$ws.Range("B827").Value = "b12f29376da282e56a56ae942e4a5f02"
$ws.Range("B830").Value = "878f501c6fcfbb24100b756563e49341"
$ws.Range("B835").Value = "44a1dc031076aedec8ddf2465a2c79d5"
